$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Fri, 03 Dec 2021 10:57:49 GMT"
$ws.Range("C5").Value = "Taiwan thing after another: the Solomon Islands"
$ws.Range("D5").Value = "00:20:44"
$ws.Range("E5").Value = "https://sphinx.acast.com/theeconomistallaudio/theintelligencepodcast/taiwanthingafteranother-thesolomonislands/media.mp3"

# Row 6
$ws.Range("A4").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Thu, 02 Dec 2021 17:00:00 GMT"
$ws.Range("C6").Value = "The Economist Asks: Eric Cantor"
$ws.Range("D6").Value = "00:27:04"
$ws.Range("E6").Value = "https://sphinx.acast.com/theeconomistallaudio/theeconomistasks/theeconomistasks-ericcantor/media.mp3"

$excel.CutCopyMode = $false
